# Adds 20 new player rows (401-420) to the "Jogadores" sheet, including a
# brand new player name "João Vitor" that becomes a new shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows: Name, Pontos(skipped/blank), Vitorias, Empate, Derrotas, Gols, Partidas, GolsSofridos, CraqueDoDia, (blank), TardeDeVitoria
$rows = @(
    @("Joãozinho",   1, 4, 2, 1, 1, 0, 0, 0, 0),
    @("Leandrinho",  1, 4, 2, 2, 1, 0, 0, 0, 0),
    @("Tom",         1, 4, 2, 0, 1, 0, 0, 0, 0),
    @("Boneco",      1, 4, 2, 1, 1, 0, 0, 0, 0),
    @("Eduardo",     1, 4, 2, 1, 1, 0, 0, 0, 0),
    @("Adriano",     4, 3, 1, 5, 1, 1, 0, 0, 0),
    @("Cabeleira",   4, 3, 1, 2, 1, 1, 0, 0, 0),
    @("Ismael",      4, 3, 1, 0, 1, 1, 0, 0, 0),
    @("Athos",       4, 3, 1, 5, 1, 1, 0, 1, 0),
    @("Corinthiano", 4, 3, 1, 1, 1, 1, 0, 0, 0),
    @("João Vitor",  1, 1, 4, 3, 1, 0, 1, 0, 0),
    @("Leandrão",    1, 1, 4, 2, 1, 0, 1, 0, 0),
    @("Marcelão",    1, 1, 4, 1, 1, 0, 1, 0, 0),
    @("Guinha",      1, 1, 4, 0, 1, 0, 1, 0, 0),
    @("Juscielio",   1, 1, 4, 1, 1, 0, 1, 0, 0),
    @("Michel",      3, 4, 1, 0, 1, 0, 0, 0, 0),
    @("Romario",     3, 4, 1, 2, 1, 0, 0, 0, 0),
    @("Miqueias",    3, 4, 1, 2, 1, 0, 0, 0, 0),
    @("Fabinho",     3, 4, 1, 3, 1, 0, 0, 0, 0),
    @("Digão",       3, 4, 1, 1, 1, 0, 0, 0, 0)
)

$startRow = 401
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]      # A - name
    $ws.Cells.Item($r, 3).Value = $data[1]      # C - Vitorias
    $ws.Cells.Item($r, 4).Value = $data[2]      # D - Empate
    $ws.Cells.Item($r, 5).Value = $data[3]      # E - Derrotas
    $ws.Cells.Item($r, 6).Value = $data[4]      # F - Gols
    $ws.Cells.Item($r, 7).Value = $data[5]      # G - Partidas
    $ws.Cells.Item($r, 8).Value = $data[6]      # H - Gols Sofridos
    $ws.Cells.Item($r, 9).Value = $data[7]      # I - Craque do Dia
    $ws.Cells.Item($r, 10).Value = $data[8]     # J
    $ws.Cells.Item($r, 11).Value = $data[9]     # K
}

# Update the sheet view to match Excel's end-state after scrolling to bottom
# of the newly added rows (selection moves to D421, directly below the data).
$ws.Range("D421").Select()
